# Adds a new "2022-Q4" quarterly sheet (right after the "总计" summary sheet
# and right before "2022-Q3") with its fund-holdings data, and records the
# corresponding new row in the "总计" summary sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q4" worksheet by copying the "2022-Q3" worksheet's
#    layout/format (same column headers / styles), inserted immediately
#    before it so the tab order becomes: 总计, 2022-Q4, 2022-Q3, 2022-Q2, ...
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3, $null)
$newSheet = $wb.ActiveSheet
$newSheet.Name = "2022-Q4"

# Clear any leftover tab-selection highlighting copied from the source sheet;
# selection itself doesn't matter much, but keep it tidy.
$newSheet.Range("A1").Select()

# ---------------------------------------------------------------------------
# 2. Populate the fund-holdings table for 2022-Q4.
#    Columns: A=index, B=基金代码, C=基金名称, D=基金规模, E=股票总仓位,
#             F=仓位占比, G=持有市值(亿元), H=仓位排名
#    B (fund codes, e.g. "009394") and D:G (percentage/amount figures, e.g.
#    "93.80") are stored as *text* in this workbook (not numbers) so leading
#    / trailing zeros survive -- force those ranges to Text format before
#    writing the values so Excel doesn't silently coerce them to numbers.
#
#    Row 8 is new (the source "2022-Q3" sheet only had 7 rows), so first
#    clone row 7's formatting down into row 8.
# ---------------------------------------------------------------------------
$newSheet.Range("A7:H7").Copy()
$newSheet.Range("A8:H8").PasteSpecial(-4122)

$newSheet.Range("B2:B8").NumberFormat = "@"
$newSheet.Range("D2:G8").NumberFormat = "@"

$rows = @(
    @(0, "161810", "银华内需精选混合（LOF）",            "23.02", "93.80", "8.46", "1.9475", 4),
    @(1, "009394", "银华同力精选混合",                    "18.05", "94.33", "6.68", "1.2057", 7),
    @(2, "320006", "诺安灵活配置混合",                    "8.35",  "77.61", "6.57", "0.5486", 3),
    @(3, "002207", "前海开源金银珠宝主题精选混合C",      "6.89",  "90.36", "7.41", "0.5105", 10),
    @(4, "001302", "前海开源金银珠宝主题精选混合A",      "3.55",  "90.36", "7.41", "0.2631", 10),
    @(5, "003304", "前海开源沪港深核心资源灵活配置混合A", "3.30",  "90.48", "7.44", "0.2455", 8),
    @(6, "003305", "前海开源沪港深核心资源灵活配置混合C", "3.17",  "90.48", "7.44", "0.2358", 8)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Value = $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
}

# ---------------------------------------------------------------------------
# 3. Add the corresponding summary row to the "总计" sheet, shifting the
#    existing quarters down by one row.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows(2).Insert()

# The newly-inserted row's formatting can differ slightly from the rest of
# the table; copy the (correct) formatting back from the row just below it.
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 7
$total.Cells.Item(2, 4).Value = 4.96

$total.Range("A1").Select()
